$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "2024-07-22 Monday" "2024-07-23 Tuesday"

Replace-Text "674÷5=" "470÷5="
Replace-Text "413÷5=" "662÷5="
Replace-Text "667÷7=" "654÷6="
Replace-Text "743÷8=" "121÷9="
Replace-Text "798÷6=" "351÷7="

Replace-Text "209÷2=" "448÷6="
Replace-Text "349÷9=" "523÷3="
Replace-Text "935÷6=" "460÷7="
Replace-Text "155÷6=" "430÷8="
Replace-Text "658÷3=" "750÷9="

Replace-Text "103÷3=" "911÷2="
Replace-Text "342÷8=" "532÷3="
Replace-Text "864÷5=" "578÷6="
Replace-Text "379÷8=" "842÷6="
Replace-Text "284÷3=" "980÷6="

Replace-Text "787÷8=" "917÷6="
Replace-Text "983÷6=" "479÷2="
Replace-Text "806÷8=" "708÷9="
Replace-Text "870÷8=" "375÷5="
Replace-Text "701÷5=" "894÷8="

Replace-Text "492÷6=" "183÷4="
Replace-Text "188÷9=" "878÷3="
Replace-Text "198÷6=" "467÷3="
Replace-Text "583÷4=" "573÷9="
Replace-Text "438÷4=" "600÷6="
